$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the indicator titles from "10.b.1" to "10.b.1.1" (Kyrgyz in A1, English in C1)
$a1 = $ws.Range("A1").Value2
$ws.Range("A1").Value = $a1.Replace("10.b.1 ", "10.b.1.1 ")

$c1 = $ws.Range("C1").Value2
$ws.Range("C1").Value = $c1.Replace("10.b.1 ", "10.b.1.1 ")

# Move/record the active cell selection as it ended up after editing
$null = $ws.Range("L8").Select()
